$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "D'Acunto, Hoang, Weber 2020 WP"
$ws.Range("B5").Value = "With a diff-in-diff, they show that households' expectations react to unconventional fiscal pol (telling them about higher taxes in the future), but not to forward guidance. The difference is due to what mechanism they understand."

# Row 6
$ws.Range("A6").Value = "Coibion, Goro, Weber 2020 WP"
$ws.Range("B6").Value = "Does policy communication during Covid influence expectations and spending plans of households? In short, no."
$ws.Range("C6").Value = "Interesting for me b/c it suggests that expectations can't incorporate communication."

# Row 7
$ws.Range("A7").Value = "Bottan, Perez-Truglia 2020 WP"
$ws.Range("B7").Value = "HHs do incorporate info about home prices and it strongly affects their decision when to sell"
$ws.Range("C7").Value = "expectations can incorporate info about individual stuff -> Preston is right, HHs understand their idiosyncratic circumstances, but not necessarily the aggregate model"

# Row 8 (shared-string build order matters: "horserace" text must be
# interned before "provide a new measure" text to reproduce target indices)
$ws.Range("A8").Value = "Bianchi, Ludvigson, Mai 2020 WP"
$ws.Range("C8").Value = "A horserace of theories seems to suggests that individuals fluctuate between optimism and pessimism (ie over- or undershooting) again, resembling learning"
$ws.Range("B8").Value = "provide a new measure of expectational errors in survey responses. Nice lit review. Interesting: their benchmark isn't RE: it's a machine-learning forecasting algorithm. Already this speaks volumes as to how reasonable learning is."

# Match styles of existing data rows (wrap text style index 2 for B/C columns)
$ws.Range("B5").WrapText = $true
$ws.Range("B6:C6").WrapText = $true
$ws.Range("B7:C7").WrapText = $true
$ws.Range("B8:C8").WrapText = $true

# Row heights: rows 5, 7, 8 are taller (wrapped multi-line text), row 6 default
$ws.Range("A5:H5").RowHeight = 30
$ws.Range("A7:H7").RowHeight = 30
$ws.Range("A8:H8").RowHeight = 30

$ws.Range("B9").Select()
